# Applies the crypto price/volume update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.617.42'
$ws.Range('E2').Value = '  +3.05%  '
$ws.Range('D3').Value = '1.850.39'
$ws.Range('E3').Value = '  +2.31%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.032'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +3.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '321.03'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.38%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.028'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.68%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4377'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.75%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3744'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.90%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07402'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.88%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8766'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.85%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '21.49'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.26%  '
$ws.Range('D12').Value = '1.861.61'
$ws.Range('E12').Value = '  -3.33%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.506'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.20%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.686'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.53%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.07170'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.18%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '82.89'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.24%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.033'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.18%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000009022'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.31%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.027'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.60%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.44'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.76%  '
$ws.Range('D21').Value = '27.637.58'
$ws.Range('E21').Value = '  +3.00%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.255'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.14%  '
$ws.Range('E23').Value = '  +0.91%  '
$ws.Range('D24').Value = '2.069.47'
$ws.Range('E24').Value = '  -3.84%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '157.34'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.02%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.935'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.89%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.74'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.71%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.293'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.75%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.949'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.68%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '116.29'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.06%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09080'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.77%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.210'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.75%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7682'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.92%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.516'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.48%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.879'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.84%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.029'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.58%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.149'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.26%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01981'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.31%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05270'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.54%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5177'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.01%  '
$ws.Range('B41').Value = 'MXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.812'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +6.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1674'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.75%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.710'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.16%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.546'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.22%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '109.04'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.74%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.57'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.83%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.720'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.19%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4657'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.03%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.06386'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.92%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.885'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.12%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '39.58'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +6.14%  '
